# Refresh the cryptocurrency price/volume snapshot (columns D and E)
# produced by the "Updated cryptos list ... with GitHub Actions" job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin (Price -> 64.094.75, Volume(1h) -> -1.03%)
$ws.Range("D2").Value = "64.094.75"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3: Ethereum (Price -> 3.152.72, Volume(1h) -> -0.38%)
$ws.Range("D3").Value = "3.152.72"
$ws.Range("E3").Value = "  -0.38%  "

# Row 4: TetherUSD (Volume(1h) -> -0.05%)
$ws.Range("E4").Value = "  -0.05%  "

# Row 5: BNB (Price -> 603.27, Volume(1h) -> -2.05%)
$ws.Range("D5").Value = "'603.27"
$ws.Range("E5").Value = "  -2.05%  "

# Row 6: Solana (Price -> 143.60, Volume(1h) -> -1.65%)
$ws.Range("D6").Value = "'143.60"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7: USDC (Volume(1h) -> -0.06%)
$ws.Range("E7").Value = "  -0.06%  "

# Row 8: LidoStakedEther (Price -> 3.146.37, Volume(1h) -> -0.45%)
$ws.Range("D8").Value = "3.146.37"
$ws.Range("E8").Value = "  -0.45%  "

# Row 9: XRP (Price -> 0.529, Volume(1h) -> -0.25%)
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -0.25%  "

# Row 10: Dogecoin (Volume(1h) -> -1.50%)
$ws.Range("E10").Value = "  -1.50%  "

# Row 11: Toncoin (Price -> 5.40, Volume(1h) -> -2.28%)
$ws.Range("D11").Value = "'5.40"
$ws.Range("E11").Value = "  -2.28%  "

# Row 12: Cardano (Price -> 0.469, Volume(1h) -> -1.08%)
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  -1.08%  "

# Row 13: ShibaInu (Price -> 0.0000255, Volume(1h) -> -1.78%)
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -1.78%  "

# Row 14: Avalanche (Price -> 35.13, Volume(1h) -> -2.02%)
$ws.Range("D14").Value = "'35.13"
$ws.Range("E14").Value = "  -2.02%  "

# Row 15: WrappedliquidstakedEther2.0 (Price -> 3.661.27, Volume(1h) -> -0.62%)
$ws.Range("D15").Value = "3.661.27"
$ws.Range("E15").Value = "  -0.62%  "

# Row 16: TRON (Volume(1h) -> +2.51%)
$ws.Range("E16").Value = "  +2.51%  "

# Row 17: WrappedBTC (Price -> 64.044.94, Volume(1h) -> -0.99%)
$ws.Range("D17").Value = "64.044.94"
$ws.Range("E17").Value = "  -0.99%  "

# Row 18: WrappedEther (Price -> 3.145.80, Volume(1h) -> -0.50%)
$ws.Range("D18").Value = "3.145.80"
$ws.Range("E18").Value = "  -0.50%  "

# Row 19: Polkadot (Price -> 6.89, Volume(1h) -> -0.51%)
$ws.Range("D19").Value = "'6.89"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20: BitcoinCash (Price -> 490.31, Volume(1h) -> +2.33%)
$ws.Range("D20").Value = "'490.31"
$ws.Range("E20").Value = "  +2.33%  "

# Row 21: Chainlink (Price -> 14.74, Volume(1h) -> +0.13%)
$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  +0.13%  "

# Row 22: Polygon (Price -> 0.715, Volume(1h) -> -0.73%)
$ws.Range("D22").Value = "'0.715"
$ws.Range("E22").Value = "  -0.73%  "

# Row 23: Uniswap (Price -> 7.78, Volume(1h) -> -2.43%)
$ws.Range("D23").Value = "'7.78"
$ws.Range("E23").Value = "  -2.43%  "

# Row 24: Litecoin (Price -> 88.35, Volume(1h) -> +4.18%)
$ws.Range("D24").Value = "'88.35"
$ws.Range("E24").Value = "  +4.18%  "

# Row 25: InternetComputer(DFINITY) (Price -> 13.35, Volume(1h) -> -3.23%)
$ws.Range("D25").Value = "'13.35"
$ws.Range("E25").Value = "  -3.23%  "

# Row 26: Dai (Volume(1h) -> +0.16%)
$ws.Range("E26").Value = "  +0.16%  "

# Row 27: PancakeSwap (Volume(1h) -> -2.04%)
$ws.Range("E27").Value = "  -2.04%  "

# Row 28: RenderToken (Price -> 8.25, Volume(1h) -> -4.20%)
$ws.Range("D28").Value = "'8.25"
$ws.Range("E28").Value = "  -4.20%  "

# Row 29: NEARProtocol (Volume(1h) -> +1.96%)
$ws.Range("E29").Value = "  +1.96%  "

# Row 30: ImmutableX (Price -> 2.08, Volume(1h) -> -0.96%)
$ws.Range("D30").Value = "'2.08"
$ws.Range("E30").Value = "  -0.96%  "

# Row 31: EthereumClassic (Price -> 27.78, Volume(1h) -> +4.20%)
$ws.Range("D31").Value = "'27.78"
$ws.Range("E31").Value = "  +4.20%  "

# Row 32: Hedera (Volume(1h) -> -5.04%)
$ws.Range("E32").Value = "  -5.04%  "

# Row 33: FirstDigitalUSD (Volume(1h) -> -0.11%)
$ws.Range("E33").Value = "  -0.11%  "

# Row 34: Stacks (Volume(1h) -> -1.37%)
$ws.Range("E34").Value = "  -1.37%  "

# Row 35: Mantle (Price -> 1.11, Volume(1h) -> -2.78%)
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  -2.78%  "

# Row 36: Filecoin (Price -> 6.09, Volume(1h) -> +0.97%)
$ws.Range("D36").Value = "'6.09"
$ws.Range("E36").Value = "  +0.97%  "

# Row 37: OKB (Price -> 52.68, Volume(1h) -> -0.98%)
$ws.Range("D37").Value = "'52.68"
$ws.Range("E37").Value = "  -0.98%  "

# Row 38: PEPE (Price -> 0.0₃0751, Volume(1h) -> -4.86%)
$ws.Range("D38").Value = "0.0₃0751"
$ws.Range("E38").Value = "  -4.86%  "

# Row 39: dogwifhat (Volume(1h) -> -7.33%)
$ws.Range("E39").Value = "  -7.33%  "

# Row 40: VeChain (Price -> 0.0399, Volume(1h) -> -0.36%)
$ws.Range("D40").Value = "'0.0399"
$ws.Range("E40").Value = "  -0.36%  "

# Row 41: Bittensor (Price -> 434.02, Volume(1h) -> -6.87%)
$ws.Range("D41").Value = "'434.02"
$ws.Range("E41").Value = "  -6.87%  "

# Row 42: Kaspa (Volume(1h) -> -0.02%)
$ws.Range("E42").Value = "  -0.02%  "

# Row 43: Cosmos (Price -> 8.39, Volume(1h) -> -0.14%)
$ws.Range("D43").Value = "'8.39"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44: Maker (Price -> 2.934.73, Volume(1h) -> +2.99%)
$ws.Range("D44").Value = "2.934.73"
$ws.Range("E44").Value = "  +2.99%  "

# Row 45: TheGraph (Volume(1h) -> -2.50%)
$ws.Range("E45").Value = "  -2.50%  "

# Row 46: Fetch.AI (Price -> 2.21, Volume(1h) -> -5.27%)
$ws.Range("D46").Value = "'2.21"
$ws.Range("E46").Value = "  -5.27%  "

# Row 47: ThetaToken (Price -> 2.42, Volume(1h) -> -1.06%)
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  -1.06%  "

# Row 48: USDe (Volume(1h) -> -0.07%)
$ws.Range("E48").Value = "  -0.07%  "

# Row 49: InjectiveProtocol (Price -> 25.97, Volume(1h) -> -2.81%)
$ws.Range("D49").Value = "'25.97"
$ws.Range("E49").Value = "  -2.81%  "

# Row 50: Stellar (Volume(1h) -> +0.67%)
$ws.Range("E50").Value = "  +0.67%  "

# Row 51: Monero (Price -> 120.47, Volume(1h) -> -0.41%)
$ws.Range("D51").Value = "'120.47"
$ws.Range("E51").Value = "  -0.41%  "
